# Update the EC workbook: remove the JAIR ALBERTO GOMEZ GALVAN worker row,
# replace the remaining worker row's data with the new period/value for
# JUAN FELIPE GOMEZ GONZALEZ PORTO, and update the summary totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 17 (CC / 1143413996 / JUAN FELIPE GOMEZ GONZALEZ PORTO / 2507 / 100000 / 2500000).
# Deleting it shifts row 16 data stays, and everything below (rows 18-23) moves up by one.
$ws.Rows("17").Delete()

# Row 16 now holds the worker that used to belong to JAIR ALBERTO GOMEZ GALVAN;
# overwrite it with the updated worker record (JUAN FELIPE GOMEZ GONZALEZ PORTO),
# new arrears period 2508, and new base salary 2500000.
$ws.Range("C16").Value = 1143413996
$ws.Range("D16").Value = "JUAN FELIPE GOMEZ GONZALEZ PORTO"
$ws.Range("E16").Value = "2508"
$ws.Range("G16").Value = 2500000

# Update the summary "VALOR MORA" total now that only one worker remains.
$ws.Range("E11").Value = 100000

# Update "Cant. Trabajadores" (worker count) from 2 to 1.
$ws.Range("C13").Value = 1
